$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "Ashley" to her full name
$ws.Range("A15").Value = 'Ashleigh "Lee" Rinkle'

# Update role/position text for Connor Haskins
$ws.Range("C6").Value = "Lead Developer Backend/ Story Board"

# Update Ashleigh's role
$ws.Range("C15").Value = "Lead Art Director/Animator/ Story Board  "

# Update role/position text for Mario Garcia
$ws.Range("C7").Value = "Project Manager/ Developer/ Art Design"

# New row introducing the "Toolsets" section further down the sheet
$ws.Range("A19").Value = "Toolsets"

# New header cell "Team Info" in row 4 (next to Engine Types table header "Unity")
$ws.Range("A4").Value = "Team Info"

# Update the selected cell to reflect where editing left off
$ws.Range("C20").Select()
